$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing data row (row 38) down onto the
# two new rows so the new cells pick up the same percentage style used by
# the rest of the dataset block.
$ws.Range("A38:B38").Copy()
$ws.Range("A39:B40").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A39").Value = "MS-C"
$ws.Range("B39").Value = 0.0553

$ws.Range("A40").Value = "MS-M"
$ws.Range("B40").Value = 0.0274

# Move the active selection to B41, mirroring where the cursor ends up
# after entering data into B39/B40 in the real worksheet.
$ws.Range("B41").Select() | Out-Null
